$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '40.790.24'
$ws.Range('E2').Value = '  +3.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.216.36'
$ws.Range('E3').Value = '  +2.51%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.82'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.633'
$ws.Range('E6').Value = '  +1.79%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.57'
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.405'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0870'
$ws.Range('E10').Value = '  +1.03%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.543.94'
$ws.Range('E12').Value = '  +2.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '16.05'
$ws.Range('E13').Value = '  +0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.22'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.822'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.63'
$ws.Range('E16').Value = '  +1.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.212.37'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '40.714.34'
$ws.Range('E18').Value = '  +3.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '73.97'
$ws.Range('E19').Value = '  +2.85%  '
$ws.Range('E20').Value = '  +5.40%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '251.83'
$ws.Range('E22').Value = '  +8.73%  '
$ws.Range('E24').Value = '  +1.19%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.31'
$ws.Range('E25').Value = '  -8.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.73'
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '173.06'
$ws.Range('E27').Value = '  +0.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.144'
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.34'
$ws.Range('E29').Value = '  +2.07%  '
$ws.Range('E30').Value = '  +2.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.83'
$ws.Range('E31').Value = '  +5.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.124'
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.67'
$ws.Range('E33').Value = '  +0.99%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.78'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('B35').Value = 'THORChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.15'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0633'
$ws.Range('E36').Value = '  +2.23%  '
$ws.Range('E37').Value = '  +5.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.47'
$ws.Range('E38').Value = '  +2.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.92'
$ws.Range('E40').Value = '  +14.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0235'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.69'
$ws.Range('E42').Value = '  +10.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.60'
$ws.Range('E43').Value = '  -2.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.24'
$ws.Range('E44').Value = '  +5.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.43'
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.508.91'
$ws.Range('E46').Value = '  -1.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.88'
$ws.Range('E47').Value = '  +2.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0940'
$ws.Range('E48').Value = '  +1.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.12'
$ws.Range('E49').Value = '  +1.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000208'
$ws.Range('E50').Value = '  +38.80%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.23'
$ws.Range('E51').Value = '  +10.49%  '
